# "new sentitivity and calculus"
# Update Model Accuracy sheet with new sensitivity columns (Market threshold,
# Market min, Market max, Recall, Precision) and refreshed accuracy numbers,
# plus refreshed confusion-matrix counts on the per-instrument sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Model Accuracy (-0.7, 0.7, 0.7)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Model Accuracy (-0.7, 0.7, 0.7)")

# New header cells, matching the existing header style (bold, bordered,
# centered) used by B1.
$ws1.Range("B1").Copy()
$ws1.Range("C1:G1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("C1").Value = "Market threshold"
$ws1.Range("D1").Value = "Market min"
$ws1.Range("E1").Value = "Market max"
$ws1.Range("F1").Value = "Recall"
$ws1.Range("G1").Value = "Precision"

# Row 2 - TOTALENERGIES SE
$ws1.Range("B2").Value = 63.08068459657702
$ws1.Range("C2").Value = 0.05450546436368681
$ws1.Range("D2").Value = -15.55441
$ws1.Range("E2").Value = 15.06418
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 0

# Row 3 - FMC CORP
$ws1.Range("B3").Value = 37.83618581907091
$ws1.Range("C3").Value = 0.009583939973006913
$ws1.Range("D3").Value = -19.35264
$ws1.Range("E3").Value = 13.70093
$ws1.Range("F3").Value = 2.144772117962467
$ws1.Range("G3").Value = 25

# Row 4 - BP PLC
$ws1.Range("B4").Value = 92.54278728606357
$ws1.Range("C4").Value = 0.04158117063764853
$ws1.Range("D4").Value = -18.75314
$ws1.Range("E4").Value = 23.33066
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 0

# Row 5 - STORA ENSO
$ws1.Range("B5").Value = 81.84596577017115
$ws1.Range("C5").Value = 0.02983403801513819
$ws1.Range("D5").Value = -12.78028
$ws1.Range("E5").Value = 12.42348
$ws1.Range("F5").Value = 0
$ws1.Range("G5").Value = 0

# Row 6 - BHP GROUP
$ws1.Range("B6").Value = 95.35452322738386
$ws1.Range("C6").Value = 0.08368817696170747
$ws1.Range("D6").Value = -16.47904
$ws1.Range("E6").Value = 14.94325
$ws1.Range("F6").Value = 0
$ws1.Range("G6").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: Confusion Matrix TOTALENERGIES SE (-0.7, 0.7, 0.7)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.7, 0.7, 0.7)")
$ws2.Range("B3").Value = 9
$ws2.Range("C3").Value = 1031
$ws2.Range("D3").Value = 9

# ---------------------------------------------------------------------------
# Sheet 3: Confusion Matrix FMC CORP (-0.7, 0.7, 0.7)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.7, 0.7, 0.7)")
$ws3.Range("B2").Value = 8
$ws3.Range("C2").Value = 18
$ws3.Range("D2").Value = 6

$ws3.Range("B3").Value = 334
$ws3.Range("C3").Value = 579
$ws3.Range("D3").Value = 318

$ws3.Range("B4").Value = 31
$ws3.Range("C4").Value = 55
$ws3.Range("D4").Value = 32

# ---------------------------------------------------------------------------
# Sheet 4: Confusion Matrix BP PLC (-0.7, 0.7, 0.7)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.7, 0.7, 0.7)")
$ws4.Range("B3").Value = 40
$ws4.Range("C3").Value = 1514
$ws4.Range("D3").Value = 42

# ---------------------------------------------------------------------------
# Sheet 5: Confusion Matrix STORA ENSO (-0.7, 0.7, 0.7)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.7, 0.7, 0.7)")
$ws5.Range("B3").Value = 109
$ws5.Range("C3").Value = 1339
$ws5.Range("D3").Value = 107

$ws5.Range("B4").Value = 1
$ws5.Range("C4").Value = 14

# ---------------------------------------------------------------------------
# Sheet 6: Confusion Matrix BHP GROUP (-0.7, 0.7, 0.7)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.7, 0.7, 0.7)")
$ws6.Range("B3").Value = 4
$ws6.Range("C3").Value = 1560
$ws6.Range("D3").Value = 3
